$d = $word.ActiveDocument

# 1) "Established AI Steering Committee" -> "Established an AI Steering Committee"
$d.Content.Find.Execute("Established AI Steering Committee", $true, $false, $false, $false, $false, $true, 1, $false, "Established an AI Steering Committee", 2)

# 2) "prioritize specific projects, " -> "prioritize projects, " (drop "specific ")
$d.Content.Find.Execute("prioritize specific projects, ", $true, $false, $false, $false, $false, $true, 1, $false, "prioritize projects, ", 2)

# 3) "saving up to 90% review time (MVP developed)" -> "saving up to 90% of review time (MVP developed)."
$d.Content.Find.Execute("saving up to 90% review time (MVP developed)", $true, $false, $false, $false, $false, $true, 1, $false, "saving up to 90% of review time (MVP developed).", 2)

# 4) "and first RAG" -> "and first retrieval-augmented generation"
$d.Content.Find.Execute("and first RAG", $true, $false, $false, $false, $false, $true, 1, $false, "and first retrieval-augmented generation", 2)

# 5) Rework the sentence about the chain-of-thought function-calling framework.
$d.Content.Find.Execute(" using the first chain-of-thought function-calling framework at the Department, enabling the tool to retrieve information from multiple source repositories and conduct variable number searches based on user queries.", $true, $false, $false, $false, $false, $true, 1, $false, " –the Department’s first chain-of-thought function-calling implementation, enabling the tool to retrieve information from multiple sources and conduct variable number searches based on user queries.", 2)

# 6) "GenAI work " -> "generative AI work " (widen match so the spellStart/spellEnd
#    proofErr markers bracketing "GenAI" get swallowed along with the edit)
$d.Content.Find.Execute("future GenAI work ", $true, $false, $false, $false, $false, $true, 1, $false, "future generative AI work ", 2)
